$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sub-total / total formulas that were previously blank cells.
$ws.Range("F8").Formula  = "=SUM(F2:F7)"
$ws.Range("F14").Formula = "=SUM(F9:F13)"
$ws.Range("F15").Formula = "=SUM(F14,F8)"
$ws.Range("F22").Formula = "=SUM(F17:F21)"
$ws.Range("F31").Formula = "=SUM(F2:F15)+SUM(F17:F22)+SUM(F26:F30)+F2"

# F32 was an empty, styled placeholder cell that is removed entirely.
$ws.Range("F32").Clear()

# Move the view / selection the way the author left it.
$ws.Range("F32").Select()
